# AFA 2020.xlsx - enter match results for rows 59-63 (RCB vs DC fixtures)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59
$ws.Range("E59").Value = 0
$ws.Range("H59").Value = 60
$ws.Range("K59").Value = 60
$ws.Range("N59").Value = 100
$ws.Range("Q59").Value = 80
$ws.Range("T59").Value = 20

# Row 60
$ws.Range("E60").Value = 20
$ws.Range("H60").Value = 80
$ws.Range("K60").Value = 60
$ws.Range("N60").Value = 100
$ws.Range("Q60").Value = 40
$ws.Range("T60").Value = 0

# Row 61
$ws.Range("E61").Value = 40
$ws.Range("H61").Value = 20
$ws.Range("K61").Value = 60
$ws.Range("N61").Value = 20
$ws.Range("Q61").Value = 100
$ws.Range("T61").Value = 80

# Row 62
$ws.Range("E62").Value = 20
$ws.Range("H62").Value = 80
$ws.Range("K62").Value = 40
$ws.Range("N62").Value = 100
$ws.Range("Q62").Value = 60
$ws.Range("T62").Value = 0

# Row 63
$ws.Range("E63").Value = 0
$ws.Range("H63").Value = 40
$ws.Range("K63").Value = 60
$ws.Range("N63").Value = 80
$ws.Range("Q63").Value = 20
$ws.Range("T63").Value = 100

$excel.Calculate()

# Ties in RANK() are not split by Excel's plain RANK formula (it gives the
# same - best - rank to every tied value), but the sheet's established
# convention (see rows 22/33/56/57) is to replace the rank-lookup formula
# with the manually averaged points for tied scores. Reproduce that here.

# Row 59: H59 and K59 are tied at 60 -> average of rank 3 & 4 points (-10,-15)
$ws.Range("G59").Value = -12.5
$ws.Range("J59").Value = -12.5

# Row 61: N61 and H61 are tied at 20 -> average of rank 5 & 6 points (-20,-25)
$ws.Range("G61").Value = -22.5
$ws.Range("M61").Value = -22.5

$excel.ActiveWindow.ScrollRow = 28
